# Update job titles/companies in the PROFESSIONAL EXPERIENCE section of the
# resume, per the new "modern resume template" described in the commit
# message. Seven Heading3 lines change their role/company text while
# keeping the location and date range intact.

$d = $word.ActiveDocument

function Replace-ExactText($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-ExactText "PRINCIPAL RESEARCH CONSULTANT - Clarity and Rigour, Washington, DC | 2012 – 2014" "DATA PRODUCTS MANAGER - Helm/Murmuration, Washington, DC | 2012 – 2014"

Replace-ExactText "DIRECTOR OF RESEARCH - Helm, Washington, DC | 2010 – 2012" "SOFTWARE ENGINEER - Mautinoa Technologies, Washington, DC | 2010 – 2012"

Replace-ExactText "SENIOR RESEARCH ANALYST - GSD&M, Austin, TX | 2008 – 2010" "SENIOR ANALYST - Myers Research, Washington, DC | 2008 – 2010"

Replace-ExactText "RESEARCH COORDINATOR - Salsa Labs, Inc., Washington, DC | 2004 – 2006" "SOFTWARE ENGINEER - Salsa Labs, Inc., Washington, DC | 2004 – 2006"

Replace-ExactText "RESEARCH MANAGER - The Praxis Project, Oakland, CA | 2002 – 2004" "INTERIM TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 – 2004"

Replace-ExactText "RESEARCH ANALYST - Lake Research Partners, Washington, DC | 2001 – 2002" "PROGRAMMER - Lake Research Partners, Washington, DC | 2001 – 2002"

Replace-ExactText "FIELD RESEARCH COORDINATOR - The Feldman Group, Washington, DC | 2000 – 2001" "FIELD DIRECTOR - The Feldman Group, Washington, DC | 2000 – 2001"

Write-Host "All replacements complete"
